$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append 5 new time log rows (15-19) for "Week 3 Work" ---
$dates  = @(42738, 42739, 42740, 42741, 42742)
$hours  = @(2, 2, 3, 2, 1)
$task   = "Week 3 Work"

# Copy the number formatting (date style) from an existing date cell (A4)
# onto the new date cells so they share the same cell style ("m/d/yyyy").
$ws.Range("A4").Copy()
$ws.Range("A15:A19").PasteSpecial(-4122)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 15 + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $task
    $ws.Cells.Item($row, 3).Value = $hours[$i]
}

# Match the selection left behind in the saved workbook
$ws.Range("C17").Select()
